# Update the revenue figures on the active sheet (row 2, columns A:M)
# to reflect the corrected year & quarter numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:M2").NumberFormat = "@"

$ws.Range("A2").Value = "€14162.725"
$ws.Range("B2").Value = "€715.28"
$ws.Range("C2").Value = "€1275.77"
$ws.Range("D2").Value = "€1089.03"
$ws.Range("E2").Value = "€738.71"
$ws.Range("F2").Value = "€1227.9"
$ws.Range("G2").Value = "€522.32"
$ws.Range("H2").Value = "€1101.1"
$ws.Range("I2").Value = "€1720.3"
$ws.Range("J2").Value = "€1614.23"
$ws.Range("K2").Value = "€1150.56"
$ws.Range("L2").Value = "€1976.16"
$ws.Range("M2").Value = "€1031.36"
